$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 607
$ws1.Range("F3").Value = 239
$ws1.Range("F8").Value = 919
$ws1.Range("F9").Value = 3824
$ws1.Range("F10").Value = 75

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 49

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 607
$ws4.Range("F3").Value = 239
$ws4.Range("F8").Value = 919
$ws4.Range("F9").Value = 3824
$ws4.Range("F10").Value = 75
$ws4.Range("F11").Value = 49
